## "new feature for date" -- rename the FECHA placeholder to DATE and
## refresh the Normal style's paragraph defaults (hyphenation / East
## Asian line-break / punctuation-overflow / spacing / alignment).

$d = $word.ActiveDocument

# --- 1. FECHA -> DATE -------------------------------------------------
# Case-sensitive match so the literal placeholder "FECHA" is replaced
# without touching the "Fecha:" label above it.
$d.Content.Find.Execute("FECHA", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "DATE", 2)

# --- 2. Normal style paragraph-format refresh -------------------------
$normal = $d.Styles("Normal")
$pf = $normal.ParagraphFormat

# Hyphenation off  -> <w:suppressAutoHyphens w:val="true"/>
$pf.Hyphenation = $false

# East-Asian line-break control off -> drops the old <w:kinsoku/>
$pf.FarEastLineBreakControl = $false

# Punctuation no longer allowed to overflow the margin
# -> <w:overflowPunct w:val="false"/> (was "true")
$pf.HangingPunctuation = $false

# Auto space between Far East and alpha text off -> drops <w:autoSpaceDE/>
$pf.AddSpaceBetweenFarEastAndAlpha = $false

# Explicit zero spacing before/after the paragraph
$pf.SpaceBefore = 0
$pf.SpaceAfter = 0

# Paragraph alignment -> <w:jc w:val="start"/>
$pf.Alignment = "start"
